$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF"), copying the existing header
# style from H1 (bold, centered, bordered) so the new header cells match
# the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the I and J data columns for rows 2-68.
$data = @(
    @(2,2,2),
    @(3,7,7),
    @(4,7,7),
    @(5,6,6),
    @(6,3,4),
    @(7,8,8),
    @(8,4,4),
    @(9,2,2),
    @(10,6,7),
    @(11,6,7),
    @(12,1,1),
    @(13,6,7),
    @(14,6,6),
    @(15,8,8),
    @(16,7,7),
    @(17,9,9),
    @(18,6,6),
    @(19,5,6),
    @(20,9,9),
    @(21,6,7),
    @(22,7,7),
    @(23,6,6),
    @(24,8,8),
    @(25,4,5),
    @(26,8,8),
    @(27,6,7),
    @(28,8,8),
    @(29,9,9),
    @(30,8,8),
    @(31,8,8),
    @(32,6,7),
    @(33,6,7),
    @(34,8,8),
    @(35,7,7),
    @(36,6,6),
    @(37,7,7),
    @(38,6,6),
    @(39,1,1),
    @(40,5,5),
    @(41,9,9),
    @(42,5,6),
    @(43,7,7),
    @(44,7,7),
    @(45,5,5),
    @(46,5,6),
    @(47,9,10),
    @(48,5,5),
    @(49,6,8),
    @(50,9,9),
    @(51,8,8),
    @(52,8,8),
    @(53,4,5),
    @(54,7,7),
    @(55,6,6),
    @(56,6,7),
    @(57,9,9),
    @(58,7,7),
    @(59,10,10),
    @(60,9,9),
    @(61,4,5),
    @(62,7,7),
    @(63,6,7),
    @(64,7,7),
    @(65,7,7),
    @(66,4,4),
    @(67,5,5),
    @(68,5,5)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 9).Value = $entry[1]
    $ws.Cells.Item($r, 10).Value = $entry[2]
}

$excel.CutCopyMode = 0
